$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row (row 1) - new columns K and L
$ws.Range("K1").Value = "fxppo2_accuracy_qkeras"
$ws.Range("L1").Value = "orig-fxppo2-drop_qkeras"

# Copy header formatting (bold, border, centered/top aligned) from existing header cell J1
$ws.Range("J1").Copy()
$ws.Range("K1:L1").PasteSpecial(-4122)

# Data values for K and L columns (rows 2-21)
$ws.Cells.Item(2, 11).Value = [double]"0.4938941655359566"
$ws.Cells.Item(2, 12).Value = [double]"-0.02170963364993217"

$ws.Cells.Item(3, 11).Value = [double]"0.004748982360922659"
$ws.Cells.Item(3, 12).Value = [double]"0.002713704206241519"

$ws.Cells.Item(4, 11).Value = [double]"0.004748982360922659"
$ws.Cells.Item(4, 12).Value = [double]"-8.673617379884035e-19"

$ws.Cells.Item(5, 11).Value = [double]"0.5040705563093623"
$ws.Cells.Item(5, 12).Value = [double]"-0.04477611940298504"

$ws.Cells.Item(6, 11).Value = [double]"0.2903663500678426"
$ws.Cells.Item(6, 12).Value = [double]"0.006784260515603768"

$ws.Cells.Item(7, 11).Value = [double]"0.004748982360922659"
$ws.Cells.Item(7, 12).Value = [double]"0.0006784260515603794"

$ws.Cells.Item(8, 11).Value = [double]"0.4477611940298508"
$ws.Cells.Item(8, 12).Value = [double]"5.551115123125783e-17"

$ws.Cells.Item(9, 11).Value = [double]"0.4477611940298508"
$ws.Cells.Item(9, 12).Value = [double]"5.551115123125783e-17"

$ws.Cells.Item(10, 11).Value = [double]"0.1024423337856174"
$ws.Cells.Item(10, 12).Value = [double]"0.06037991858887382"

$ws.Cells.Item(11, 11).Value = [double]"0.4864314789687924"
$ws.Cells.Item(11, 12).Value = [double]"-0.03459972862957938"

$ws.Cells.Item(12, 11).Value = [double]"0.2740841248303935"
$ws.Cells.Item(12, 12).Value = [double]"0.1146540027137042"

$ws.Cells.Item(13, 11).Value = [double]"0.2815468113975577"
$ws.Cells.Item(13, 12).Value = [double]"0.004070556309362261"

$ws.Cells.Item(14, 11).Value = [double]"0.4484396200814111"
$ws.Cells.Item(14, 12).Value = [double]"0.03188602442333788"

$ws.Cells.Item(15, 11).Value = [double]"0.4477611940298508"
$ws.Cells.Item(15, 12).Value = [double]"5.551115123125783e-17"

$ws.Cells.Item(16, 11).Value = [double]"0.2971506105834464"
$ws.Cells.Item(16, 12).Value = [double]"0"

$ws.Cells.Item(17, 11).Value = [double]"0.4728629579375848"
$ws.Cells.Item(17, 12).Value = [double]"-0.0237449118046133"

$ws.Cells.Item(18, 11).Value = [double]"0.4138398914518318"
$ws.Cells.Item(18, 12).Value = [double]"0.06987788331071915"

$ws.Cells.Item(19, 11).Value = [double]"0.4477611940298508"
$ws.Cells.Item(19, 12).Value = [double]"5.551115123125783e-17"

$ws.Cells.Item(20, 11).Value = [double]"0.4308005427408412"
$ws.Cells.Item(20, 12).Value = [double]"0.04206241519674359"

$ws.Cells.Item(21, 11).Value = [double]"0.4972862957937585"
$ws.Cells.Item(21, 12).Value = [double]"-0.04206241519674359"
